$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.914.44"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.419.99"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.25"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.10"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.403.73"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.02"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.827.42"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.811.16"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.402.62"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.63"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("E20").Value = "  +3.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.41"
$ws.Range("E21").Value = "  -1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.04"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.52"
$ws.Range("E27").Value = "  -7.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "578.03"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.515.79"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0922"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("E32").Value = "  -5.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -5.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.40"
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "149.69"
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.30"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.16"
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("E43").Value = "  -3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.14"
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("E45").Value = "  -5.17%  "
$ws.Range("E46").Value = "  +12.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.44"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.61"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0506"
$ws.Range("E51").Value = "  -3.31%  "
